$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45810
$ws.Range("B2").Value = 45814

$ws.Range("K3").Select()
